# Apply the fixture update: point the "derivative" file cell at the
# already-existing video fixture name, set the sheet tab color, and
# move the active selection to the cell that was edited (E3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell E3 ("File" column, row 3) held the shared string "derivative.mp4".
# Re-use the existing fixture name "videoshort.mp4" instead.
$ws.Range("E3").Value = "videoshort.mp4"

# Give the sheet tab a (white) color. COM Color values are packed as
# 0x00BBGGRR, so white (255,255,255) is simply 16777215.
$ws.Tab.Color = 16777215

# Move the selection/active cell onto the edited cell.
$ws.Range("E3").Select()
